$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts")

# Fill the "Definition" column (D) with the same value as the
# "Display" column (C) for each concept row.
for ($row = 2; $row -le 4; $row++) {
    $ws.Cells.Item($row, 4).Value = $ws.Cells.Item($row, 3).Value2
}
